# Add a new ToDo row (row 3) to the Trashcan worksheet, mirroring the
# existing row 2 for the "deadline / done / importance" columns while
# introducing a new subject/task pair (fixing a PMD rule violation per
# #1245).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the mangled-encoding task text exactly as produced by the original
# (corrupted IME) input: same garbled prefix as row 2's task cell, but
# ending in "qq" instead of ".zzzzzzzz".
$c1  = [char]0xFFFD
$c2  = [char]0x0537
$c3  = [char]0xFFFD
$c4  = [char]0xFFFD
$c5  = [char]0xFFFD
$c6  = [char]0xFFFD
$c7  = [char]0x05BC
$c8  = [char]0xFFFD
$c9  = [char]0xFFFD
$c10 = [char]0xFFFD
$taskText = "$c1$c2$c3$c4$c5$c6$c7$c8$c9$c10" + "qq"

$ws.Range("A3").Value2 = "알고리즘"
$ws.Range("B3").Value2 = $taskText

# Columns C-F repeat row 2's values verbatim; copy the cells instead of
# re-typing so the shared-string / numeric typing matches exactly (F2 is
# the text "3", which would otherwise be re-interpreted as a number).
$ws.Range("C2").Copy($ws.Range("C3"))
$ws.Range("D2").Copy($ws.Range("D3"))
$ws.Range("E2").Copy($ws.Range("E3"))
$ws.Range("F2").Copy($ws.Range("F3"))
